$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "C2" = 3.448290604152642
    "D2" = 9.145498974178896
    "E2" = 14.90729549833441
    "F2" = 22.69786651950371
    "G2" = 21.08608397561744
    "H2" = 12.33454275493188
    "I2" = 16.86886856337481
    "J2" = 10.15345392877687
    "M2" = 58.96029491922199
    "O2" = 17.74769338902397
    "C3" = 3.311900644384297
    "D3" = 9.208199912639609
    "E3" = 14.65333387160127
    "F3" = 23.16121294231622
    "G3" = 21.3898229122545
    "H3" = 12.45924227866217
    "I3" = 17.04477107235396
    "J3" = 10.03852328933165
    "M3" = 55.47616486653173
    "O3" = 17.97272643372167
    "C4" = 3.224649720217164
    "D4" = 9.25061306240695
    "E4" = 14.50273420706088
    "F4" = 23.45879694974228
    "G4" = 21.59711826213518
    "H4" = 12.53990921554777
    "I4" = 17.16127376091072
    "J4" = 9.971630773580079
    "M4" = 53.21265024487478
    "O4" = 18.11985177749606
    "C5" = 3.188248436766087
    "D5" = 9.268867304283075
    "E5" = 14.44277968521281
    "F5" = 23.5833350129184
    "G5" = 21.68665683490087
    "H5" = 12.57380762152682
    "I5" = 17.21085354530042
    "J5" = 9.945318465558781
    "M5" = 52.25916909854443
    "O5" = 18.18202748088657
    "C6" = 3.182154094436265
    "D6" = 9.271956632669101
    "E6" = 14.43291195026449
    "F6" = 23.60421140851266
    "G6" = 21.70182575468127
    "H6" = 12.57949824406054
    "I6" = 17.21921236580119
    "J6" = 9.941007150970988
    "M6" = 52.09897313604252
    "O6" = 18.19248497890953
    "C7" = 3.224162173321179
    "D7" = 9.250855333619384
    "E7" = 14.50191980860182
    "F7" = 23.46046329854
    "G7" = 21.5983055217509
    "H7" = 12.54036223387208
    "I7" = 17.16193393483566
    "J7" = 9.971272053655229
    "M7" = 53.1999168472071
    "O7" = 18.12068134770277
    "C8" = 3.402011349016229
    "D8" = 9.166297326945704
    "E8" = 14.81866801071761
    "F8" = 22.85489890444464
    "G8" = 21.18639505761549
    "H8" = 12.37668478859577
    "I8" = 16.9277359510538
    "J8" = 10.11307714195872
    "M8" = 57.78483120405355
    "O8" = 17.82340747399664
    "C9" = 3.721598822270434
    "D9" = 9.032186346995637
    "E9" = 15.47877087814714
    "F9" = 21.77227401486503
    "G9" = 20.55133960773465
    "H9" = 12.0884960080452
    "I9" = 16.53742127643217
    "J9" = 10.41925259367044
    "M9" = 65.78596410415629
    "O9" = 17.31294547681288
    "C10" = 3.937102013111749
    "D10" = 8.953916541184171
    "E10" = 15.98298685313068
    "F10" = 21.042452566011
    "G10" = 20.2008836962199
    "H10" = 11.8970989458039
    "I10" = 16.29485282555924
    "J10" = 10.65983042486402
    "M10" = 71.0589622681975
    "O10" = 16.98420226080889
    "C11" = 4.030703274271977
    "D11" = 8.922923242792677
    "E11" = 16.21560836313866
    "F11" = 20.72511807043771
    "G11" = 20.06907395248145
    "H11" = 11.81452661585146
    "I11" = 16.19459101573437
    "J11" = 10.77233765609649
    "M11" = 73.32664054715393
    "O11" = 16.84519523157907
    "C12" = 4.065493785185661
    "D12" = 8.911868723161641
    "E12" = 16.30408834580793
    "F12" = 20.60710522818806
    "G12" = 20.02333399162906
    "H12" = 11.78391285297322
    "I12" = 16.1581172988394
    "J12" = 10.81535398286433
    "M12" = 74.16656358688364
    "O12" = 16.79411700422427
    "C13" = 4.058030365384582
    "D13" = 8.914218874492773
    "E13" = 16.28501637394759
    "F13" = 20.63242476757003
    "G13" = 20.03299605540539
    "H13" = 11.79047682532789
    "I13" = 16.16590544929339
    "J13" = 10.80607176625939
    "M13" = 73.98650621852184
    "O13" = 16.80504745499605
    "C14" = 4.033578710187591
    "D14" = 8.922000003149561
    "E14" = 16.22288024762448
    "F14" = 20.71536561067384
    "G14" = 20.0652260905688
    "H14" = 11.81199483542937
    "I14" = 16.19156012121439
    "J14" = 10.77586853645374
    "M14" = 73.39611791881522
    "O14" = 16.84096144841273
    "C15" = 4.018515705840212
    "D15" = 8.926855550845641
    "E15" = 16.18486882562705
    "F15" = 20.76645136914134
    "G15" = 20.08551746618093
    "H15" = 11.82526073683606
    "I15" = 16.20747012258114
    "J15" = 10.75742102720649
    "M15" = 73.03204174355994
    "O15" = 16.86316444081582
    "C16" = 3.930894255680129
    "D16" = 8.956036419668319
    "E16" = 15.96784331837313
    "F16" = 21.06349001390855
    "G16" = 20.21007073503677
    "H16" = 11.90258638326877
    "I16" = 16.30161174239461
    "J16" = 10.65253706762981
    "M16" = 70.90813077442203
    "O16" = 16.99350258223358
    "C17" = 3.875993233328426
    "D17" = 8.975132312912342
    "E17" = 15.83548262057634
    "F17" = 21.24950049586836
    "G17" = 20.29369652695613
    "H17" = 11.951180134392
    "I17" = 16.36197755258038
    "J17" = 10.58895920591247
    "M17" = 69.57164421227009
    "O17" = 17.07619159879498
    "C18" = 3.843999261467098
    "D18" = 8.986548704473078
    "E18" = 15.75966263864309
    "F18" = 21.35786603334093
    "G18" = 20.34438281140936
    "H18" = 11.97955255787388
    "I18" = 16.39764482673924
    "J18" = 10.55268131914382
    "M18" = 68.7905946390124
    "O18" = 17.12474134648136
    "C19" = 3.833095702731521
    "D19" = 8.990487896109068
    "E19" = 15.73404698207805
    "F19" = 21.39479197872626
    "G19" = 20.36198247582862
    "H19" = 11.98923131136232
    "I19" = 16.40988241547592
    "J19" = 10.54044900926187
    "M19" = 68.52402461236279
    "O19" = 17.1413480233939
    "C20" = 3.881880763163185
    "D20" = 8.973054581316264
    "E20" = 15.84954108150246
    "F20" = 21.22955651070855
    "G20" = 20.2845252165558
    "H20" = 11.94596345332161
    "I20" = 16.35545324610652
    "J20" = 10.59569732983023
    "M20" = 69.71519214852641
    "O20" = 17.06728651377152
    "C21" = 4.040778631878604
    "D21" = 8.919695835397834
    "E21" = 16.24112110579497
    "F21" = 20.69094499446087
    "G21" = 20.05564445970918
    "H21" = 11.80565664092113
    "I21" = 16.18398384375683
    "J21" = 10.78472899600885
    "M21" = 73.57003897041598
    "O21" = 16.83036990808785
    "C22" = 4.140808043735234
    "D22" = 8.888807070486955
    "E22" = 16.49928908764409
    "F22" = 20.35151733036624
    "G22" = 19.93047854069522
    "H22" = 11.71777766374072
    "I22" = 16.08064363013417
    "J22" = 10.9106621048891
    "M22" = 75.97987272826221
    "O22" = 16.68465108824414
    "C23" = 4.087774863628141
    "D23" = 8.90492203883019
    "E23" = 16.36131857146363
    "F23" = 20.53150774845419
    "G23" = 19.99498031029043
    "H23" = 11.76432797867251
    "I23" = 16.13498524371143
    "J23" = 10.84323997477531
    "M23" = 74.70370113555688
    "O23" = 16.76157390857042
    "C24" = 3.879220350970489
    "D24" = 8.973992561196939
    "E24" = 15.84318438598047
    "F24" = 21.23856874927124
    "G24" = 20.2886634558059
    "H24" = 11.94832055980627
    "I24" = 16.3583998913857
    "J24" = 10.59265016876746
    "M24" = 69.65033364270067
    "O24" = 17.07130935789144
    "C25" = 3.638436208122109
    "D25" = 9.064978754121784
    "E25" = 15.2964924236714
    "F25" = 22.05376896226921
    "G25" = 20.70354914201533
    "H25" = 12.16291446112615
    "I25" = 16.63541212502446
    "J25" = 10.33356735903747
    "M25" = 63.7277733032323
    "O25" = 17.44306909787668
}

foreach ($key in $data.Keys) {
    $ws.Range($key).Value = $data[$key]
}
